# Auto-generated Excel COM-interop script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.971.66'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '2.264.44'
$ws.Range("E3").Value = '  -0.26%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = "'305.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.09%  '
$ws.Range("D6").Value = "'95.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.68%  '
$ws.Range("E7").Value = '  -0.64%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").Value = "'0.490"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.67%  '
$ws.Range("D10").Value = "'35.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.43%  '
$ws.Range("D11").Value = "'0.0789"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.09%  '
$ws.Range("E12").Value = '  -0.44%  '
$ws.Range("D13").Value = "'6.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.76%  '
$ws.Range("D14").Value = '2.612.78'
$ws.Range("E14").Value = '  -0.29%  '
$ws.Range("D15").Value = "'14.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.65%  '
$ws.Range("D16").Value = '2.255.94'
$ws.Range("E16").Value = '  -0.71%  '
$ws.Range("D17").Value = "'0.792"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.17%  '
$ws.Range("D18").Value = '41.869.43'
$ws.Range("E18").Value = '  +0.19%  '
$ws.Range("D19").Value = "'12.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.02%  '
$ws.Range("D20").Value = '0.0₃0902'
$ws.Range("E20").Value = '  -2.06%  '
$ws.Range("D21").Value = "'5.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.32%  '
$ws.Range("D22").Value = "'67.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.47%  '
$ws.Range("D23").Value = "'237.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.69%  '
$ws.Range("D24").Value = "'2.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.72%  '
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("E26").Value = '  -1.00%  '
$ws.Range("D27").Value = "'23.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.46%  '
$ws.Range("D28").Value = "'36.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.24%  '
$ws.Range("D29").Value = "'9.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.76%  '
$ws.Range("D30").Value = "'2.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.80%  '
$ws.Range("D31").Value = "'160.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.49%  '
$ws.Range("D32").Value = "'5.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.53%  '
$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("E34").Value = '  +4.97%  '
$ws.Range("D35").Value = "'0.0737"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.98%  '
$ws.Range("D36").Value = "'17.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.32%  '
$ws.Range("E37").Value = '  +0.51%  '
$ws.Range("E38").Value = '  -0.63%  '
$ws.Range("E39").Value = '  +1.71%  '
$ws.Range("E40").Value = '  -2.17%  '
$ws.Range("D41").Value = "'4.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.69%  '
$ws.Range("E42").Value = '  +6.43%  '
$ws.Range("D43").Value = '1.976.70'
$ws.Range("E43").Value = '  -1.77%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = "'18.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.37%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").Value = "'0.0283"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.16%  '
$ws.Range("D46").Value = "'2.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.78%  '
$ws.Range("D47").Value = "'9.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.73%  '
$ws.Range("D48").Value = "'53.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.81%  '
$ws.Range("D49").Value = "'72.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.39%  '
$ws.Range("D50").Value = "'1.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.22%  '
$ws.Range("D51").Value = "'90.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.21%  '
